$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'31.039.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.78%  "
$ws.Range("D3").Value = "'1.685.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.15%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'220.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.44%  "
$ws.Range("D6").Value = "'0.532"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.09%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'29.32"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.08%  "
$ws.Range("E9").Value = "  +2.41%  "
$ws.Range("D10").Value = "'0.0638"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.56%  "
$ws.Range("D11").Value = "'0.0909"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").Value = "'1.926.91"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.20%  "
$ws.Range("D13").Value = "'1.687.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.26%  "
$ws.Range("D14").Value = "'10.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.31%  "
$ws.Range("E15").Value = "  +3.64%  "
$ws.Range("D16").Value = "'4.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.50%  "
$ws.Range("D17").Value = "'31.003.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.65%  "
$ws.Range("D18").Value = "'66.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.95%  "
$ws.Range("D19").Value = "'247.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.61%  "
$ws.Range("D20").Value = "'0.0₃0721"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.37%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "'4.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.44%  "
$ws.Range("D23").Value = "'10.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.70%  "
$ws.Range("E24").Value = "  -1.07%  "
$ws.Range("D25").Value = "'158.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("E26").Value = "  +2.57%  "
$ws.Range("E27").Value = "  +2.46%  "
$ws.Range("D28").Value = "'6.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  +2.09%  "
$ws.Range("E31").Value = "  +3.86%  "
$ws.Range("D32").Value = "'3.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.97%  "
$ws.Range("D33").Value = "'3.34"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.69%  "
$ws.Range("D34").Value = "'1.513.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.22%  "
$ws.Range("E35").Value = "  +2.58%  "
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").Value = "'83.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.11%  "
$ws.Range("D38").Value = "'0.613"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.76%  "
$ws.Range("E39").Value = "  +4.84%  "
$ws.Range("D40").Value = "'2.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.62%  "
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").Value = "'2.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.842"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("D44").Value = "'0.0504"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("E45").Value = "  +2.74%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "'51.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.99%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'5.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.28%  "
$ws.Range("D49").Value = "'1.821.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'93.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.32%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "'0.0₆0121"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.31%  "
